$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the long step-by-step instruction text in column D (SCENARIO_DESC)
# with the shorter "<Action> Setup Jenis Beban" labels, and shrink the
# row heights to match (the text no longer needs to wrap across several lines).
$ws.Range("D2").Value = "Tambah Setup Jenis Beban"
$ws.Rows(2).RowHeight = 30

$ws.Range("D3").Value = "View Setup Jenis Beban"
$ws.Rows(3).RowHeight = 30

$ws.Range("D4").Value = "Ubah Setup Jenis Beban"
$ws.Rows(4).AutoFit()

$ws.Range("D5").Value = "Hapus Setup Jenis Beban"
$ws.Rows(5).RowHeight = 30

# Move the active selection to D5 (matches the saved sheet view state).
$ws.Range("D5").Select()
